$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3449.45
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 3588.8948
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 10766.6844
$ws.Range("M17").Value = -2232
$ws.Range("N17").Value = -11102.6844
$ws.Range("H32").Value = 1094.4
$ws.Range("I32").Value = 850
$ws.Range("J32").Value = 1155.5
$ws.Range("K32").Value = 850
$ws.Range("L32").Value = 1155.5
$ws.Range("M32").Value = -524
$ws.Range("N32").Value = -1807.5
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("N46").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("N60").Value = 0
$ws.Range("H129").Value = 2190.25
$ws.Range("I129").Value = 1555.8572
$ws.Range("K129").Value = 4667.571599999999
$ws.Range("M129").Value = 332.4284000000007
$ws.Range("H133").Value = 119962
$ws.Range("J133").Value = 119962
$ws.Range("L133").Value = 119962
$ws.Range("N133").Value = -130082
$ws.Range("H137").Value = 2619.4138
$ws.Range("I137").Value = 2184.0527
$ws.Range("K137").Value = 6552.158100000001
$ws.Range("M137").Value = -4002.158100000001
$ws.Range("H138").Value = 2697.647
$ws.Range("I138").Value = 2576.2222
$ws.Range("J138").Value = 2834.25
$ws.Range("K138").Value = 7728.6666
$ws.Range("L138").Value = 8502.75
$ws.Range("M138").Value = -2588.6666
$ws.Range("N138").Value = -18782.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 30189.953
$ws.Range("I32").Value = 31899.21
$ws.Range("J32").Value = 17199.6
$ws.Range("K32").Value = 31899.21
$ws.Range("L32").Value = 17199.6
$ws.Range("M32").Value = -31612.21
$ws.Range("N32").Value = -17773.6
$ws.Range("H61").Value = 5662.951
$ws.Range("I61").Value = 3520.8484
$ws.Range("J61").Value = 14499.125
$ws.Range("K61").Value = 3520.8484
$ws.Range("L61").Value = 14499.125
$ws.Range("M61").Value = -3308.8484
$ws.Range("N61").Value = -14923.125
$ws.Range("H74").Value = 3001.5217
$ws.Range("I74").Value = 1272.0834
$ws.Range("J74").Value = 4888.1816
$ws.Range("K74").Value = 1272.0834
$ws.Range("L74").Value = 4888.1816
$ws.Range("M74").Value = -398.0834
$ws.Range("N74").Value = -6636.1816
$ws.Range("H77").Value = 3001.5217
$ws.Range("I77").Value = 1272.0834
$ws.Range("J77").Value = 4888.1816
$ws.Range("K77").Value = 6360.416999999999
$ws.Range("L77").Value = 24440.908
$ws.Range("M77").Value = -1992.416999999999
$ws.Range("N77").Value = -33176.908
$ws.Range("H102").Value = 12552157
$ws.Range("I102").Value = 2478.647
$ws.Range("K102").Value = 2478.647
$ws.Range("M102").Value = -856.6469999999999
$ws.Range("H110").Value = 17859380
$ws.Range("I110").Value = 27779036
$ws.Range("J110").Value = 4000
$ws.Range("K110").Value = 27779036
$ws.Range("L110").Value = 4000
$ws.Range("M110").Value = -27776991
$ws.Range("N110").Value = -8090
$ws.Range("H132").Value = 3870.5
$ws.Range("I132").Value = 2796.147
$ws.Range("K132").Value = 8388.440999999999
$ws.Range("M132").Value = -5858.440999999999
$ws.Range("H136").Value = 5662.951
$ws.Range("I136").Value = 3520.8484
$ws.Range("J136").Value = 14499.125
$ws.Range("K136").Value = 10562.5452
$ws.Range("L136").Value = 43497.375
$ws.Range("M136").Value = -8012.5452
$ws.Range("N136").Value = -48597.375

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 75282.74000000001
$ws.Range("I86").Value = 1299.762
$ws.Range("J86").Value = 334223.16
$ws.Range("K86").Value = 1299.762
$ws.Range("L86").Value = 334223.16
$ws.Range("M86").Value = -176.7619999999999
$ws.Range("N86").Value = -336469.16
$ws.Range("H89").Value = 75282.74000000001
$ws.Range("I89").Value = 1299.762
$ws.Range("J89").Value = 334223.16
$ws.Range("K89").Value = 6498.809999999999
$ws.Range("L89").Value = 1671115.8
$ws.Range("M89").Value = -882.8099999999995
$ws.Range("N89").Value = -1682347.8
$ws.Range("H94").Value = 1532.2222
$ws.Range("I94").Value = 1455.3529
$ws.Range("J94").Value = 1662.9
$ws.Range("K94").Value = 1455.3529
$ws.Range("L94").Value = 1662.9
$ws.Range("M94").Value = -1004.3529
$ws.Range("N94").Value = -2564.9
$ws.Range("H107").Value = 1395.5
$ws.Range("I107").Value = 1320.6428
$ws.Range("J107").Value = 1657.5
$ws.Range("K107").Value = 1320.6428
$ws.Range("L107").Value = 1657.5
$ws.Range("M107").Value = 599.3571999999999
$ws.Range("N107").Value = -5497.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37041176
$ws.Range("I31").Value = 142858220
$ws.Range("J31").Value = 5209.75
$ws.Range("K31").Value = 142858220
$ws.Range("L31").Value = 5209.75
$ws.Range("M31").Value = -142857925
$ws.Range("N31").Value = -5799.75
$ws.Range("H34").Value = 37041176
$ws.Range("I34").Value = 142858220
$ws.Range("J34").Value = 5209.75
$ws.Range("K34").Value = 142858220
$ws.Range("L34").Value = 5209.75
$ws.Range("M34").Value = -142858018
$ws.Range("N34").Value = -5613.75
$ws.Range("H69").Value = 5998
$ws.Range("I69").Value = 5998
$ws.Range("K69").Value = 5998
$ws.Range("M69").Value = -5249
$ws.Range("H72").Value = 5998
$ws.Range("I72").Value = 5998
$ws.Range("K72").Value = 17994
$ws.Range("M72").Value = -14250
$ws.Range("H134").Value = 3905.9666
$ws.Range("I134").Value = 2801.9
$ws.Range("J134").Value = 6114.1
$ws.Range("K134").Value = 8405.700000000001
$ws.Range("L134").Value = 18342.3
$ws.Range("M134").Value = -5870.700000000001
$ws.Range("N134").Value = -23412.3

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4999.5
$ws.Range("I70").Value = 4999.5
$ws.Range("K70").Value = 14998.5
$ws.Range("M70").Value = -14683.5
$ws.Range("H73").Value = 4999.5
$ws.Range("I73").Value = 4999.5
$ws.Range("K73").Value = 14998.5
$ws.Range("M73").Value = -13906.5
$ws.Range("H80").Value = 3750
$ws.Range("J80").Value = 4025
$ws.Range("L80").Value = 12075
$ws.Range("N80").Value = -13947
$ws.Range("H83").Value = 3750
$ws.Range("J83").Value = 4025
$ws.Range("L83").Value = 36225
$ws.Range("N83").Value = -45585
$ws.Range("H94").Value = 5000
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352
$ws.Range("H109").Value = 2070.3333
$ws.Range("I109").Value = 1730.75
$ws.Range("K109").Value = 5192.25
$ws.Range("M109").Value = -4152.25
$ws.Range("H131").Value = 7578415
$ws.Range("I131").Value = 13889734
$ws.Range("K131").Value = 41669202
$ws.Range("M131").Value = -41664162

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 250986.08
$ws.Range("I14").Value = 542481.5
$ws.Range("J14").Value = 1132.8572
$ws.Range("K14").Value = 542481.5
$ws.Range("L14").Value = 1132.8572
$ws.Range("M14").Value = -542313.5
$ws.Range("N14").Value = -1468.8572
$ws.Range("H107").Value = 451.58823
$ws.Range("I107").Value = 489.5
$ws.Range("J107").Value = 397.42856
$ws.Range("K107").Value = 489.5
$ws.Range("L107").Value = 397.42856
$ws.Range("M107").Value = 1430.5
$ws.Range("N107").Value = -4237.42856
$ws.Range("H122").Value = 11338
$ws.Range("I122").Value = 5342.6665
$ws.Range("K122").Value = 16027.9995
$ws.Range("M122").Value = -13577.9995
$ws.Range("H132").Value = 5039.5557
$ws.Range("I132").Value = 3898.3618
$ws.Range("K132").Value = 11695.0854
$ws.Range("M132").Value = -9165.0854

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 14513
$ws.Range("J3").Value = 14513
$ws.Range("L3").Value = 14513
$ws.Range("N3").Value = -14737
$ws.Range("H15").Value = 14513
$ws.Range("J15").Value = 14513
$ws.Range("L15").Value = 14513
$ws.Range("N15").Value = -14853
$ws.Range("H21").Value = 2416.6667
$ws.Range("J21").Value = 2416.6667
$ws.Range("L21").Value = 2416.6667
$ws.Range("N21").Value = -2764.6667
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H132").Value = 3422
$ws.Range("I132").Value = 2823.0417
$ws.Range("K132").Value = 8469.125100000001
$ws.Range("M132").Value = -5939.125100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3688.8
$ws.Range("I136").Value = 1404.4166
$ws.Range("J136").Value = 5797.4614
$ws.Range("K136").Value = 4213.2498
$ws.Range("L136").Value = 17392.3842
$ws.Range("M136").Value = -1663.2498
$ws.Range("N136").Value = -22492.3842

